$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "İsim" header / name list from column A to column B,
# and add a new "Numara" (number) column in column A.

# Header row
$ws.Range("B1").Value = "İsim"
$ws.Range("A1").Value = "Numara"

# Move names from column A into column B
$ws.Range("B2").Value = "ahmet"
$ws.Range("B3").Value = "mehmet"
$ws.Range("B4").Value = "ali"
$ws.Range("B5").Value = "fatma"

# Add the new "Numara" values in column A
$ws.Range("A2").Value = 201008
$ws.Range("A3").Value = 201009
$ws.Range("A4").Value = 201010
$ws.Range("A5").Value = 201011

$ws.Range("B5").Select()
